$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear D2 (was "33,33 TL - 33,33 TL")
$ws.Range("D2").Value = ""

# Update J2 from "25 TL - 25 TL" to "65 TL - 65 TL"
$ws.Range("J2").Value = "65 TL - 65 TL"

# Clear D7 (was "%1,6")
$ws.Range("D7").Value = ""
